$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for I0 and IF columns, copying style from existing header cell (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data rows 2-74 for columns I (I0) and J (IF)
$data = @(
    @(2,5,6),
    @(3,6,6),
    @(4,4,5),
    @(5,6,6),
    @(6,6,6),
    @(7,6,6),
    @(8,6,6),
    @(9,8,8),
    @(10,7,7),
    @(11,10,10),
    @(12,6,6),
    @(13,9,9),
    @(14,6,6),
    @(15,9,9),
    @(16,9,9),
    @(17,10,10),
    @(18,7,7),
    @(19,7,8),
    @(20,7,7),
    @(21,6,6),
    @(22,9,9),
    @(23,9,9),
    @(24,7,7),
    @(25,5,5),
    @(26,6,6),
    @(27,9,9),
    @(28,8,8),
    @(29,6,6),
    @(30,9,9),
    @(31,7,8),
    @(32,7,7),
    @(33,8,8),
    @(34,8,8),
    @(35,8,8),
    @(36,8,8),
    @(37,7,7),
    @(38,7,8),
    @(39,8,8),
    @(40,6,7),
    @(41,8,8),
    @(42,8,8),
    @(43,8,8),
    @(44,8,8),
    @(45,9,9),
    @(46,8,8),
    @(47,7,7),
    @(48,8,8),
    @(49,8,8),
    @(50,7,7),
    @(51,7,8),
    @(52,7,7),
    @(53,8,8),
    @(54,7,8),
    @(55,6,6),
    @(56,9,9),
    @(57,6,7),
    @(58,8,9),
    @(59,8,8),
    @(60,9,9),
    @(61,6,7),
    @(62,8,8),
    @(63,8,8),
    @(64,7,8),
    @(65,7,8),
    @(66,4,6),
    @(67,5,5),
    @(68,9,9),
    @(69,6,7),
    @(70,8,8),
    @(71,8,9),
    @(72,9,9),
    @(73,3,3),
    @(74,2,2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

